$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4156708121299744
$ws.Range("B1").Value = 0.4032054841518402
$ws.Range("C1").Value = 0.4150528013706207
$ws.Range("D1").Value = 0.5502541661262512
$ws.Range("E1").Value = 0.708649754524231
